$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right above the current first data block row (461),
# shifting all existing data rows (461:561) down to (463:563).
$ws.Rows("461:462").Insert()

# New row 461 - "Primera" quality entry for the new reporting date.
$ws.Range("A461").Value = 1
$ws.Range("B461").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C461").Value = "Arica y Parinacota"
$ws.Range("D461").Value = 45275
$ws.Range("E461").Value = 15
$ws.Range("F461").Value = 100112043
$ws.Range("G461").Value = "Pepino ensalada"
$ws.Range("H461").Value = "Sin especificar"
$ws.Range("I461").Value = "Primera"
$ws.Range("J461").Value = 200
$ws.Range("K461").Value = 8000
$ws.Range("L461").Value = 9000
$ws.Range("M461").Value = 8500
$ws.Range("N461").Value = "$/caja 70 unidades"
$ws.Range("O461").Value = "Región de Arica y Parinacota"
$ws.Range("P461").Value = 121
$ws.Range("Q461").Value = 70
$ws.Range("R461").Value = "Hortaliza"

# New row 462 - "Segunda" quality entry for the new reporting date.
$ws.Range("A462").Value = 1
$ws.Range("B462").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C462").Value = "Arica y Parinacota"
$ws.Range("D462").Value = 45275
$ws.Range("E462").Value = 15
$ws.Range("F462").Value = 100112043
$ws.Range("G462").Value = "Pepino ensalada"
$ws.Range("H462").Value = "Sin especificar"
$ws.Range("I462").Value = "Segunda"
$ws.Range("J462").Value = 350
$ws.Range("K462").Value = 6000
$ws.Range("L462").Value = 7000
$ws.Range("M462").Value = 6429
$ws.Range("N462").Value = "$/caja 100 unidades"
$ws.Range("O462").Value = "Región de Arica y Parinacota"
$ws.Range("P462").Value = 64
$ws.Range("Q462").Value = 100
$ws.Range("R462").Value = "Hortaliza"
